$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the peripheral/pin name (ADC channel, timer channel, UART) to the
# "use/description" text in column D for the rows that changed when the
# adc_sample project moved to the cpp framework.
$ws.Range("D2").Value = "雷达中频输入ADC2_IN5"
$ws.Range("D3").Value = "环境光输入ADC1_IN3"
$ws.Range("D7").Value = "FSK控制输出TIM1_CH3"
$ws.Range("D9").Value = "串口发uart3"
$ws.Range("D10").Value = "串口收uart3"

# Match the saved cursor position recorded in the file.
[void]$ws.Range("D13").Select()
